$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("DesignNo","Jeweltype","Category","Subcategory","Item","Procatgory","weight","Purity","color","size","style","unit","making","qty")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
}

$headerRange = $ws.Range("A1:N1")
$headerRange.Font.Bold = $true
$headerRange.Font.Size = 9.85
$headerRange.Font.Name = "Arial"
$headerRange.Font.Color = 0
$headerRange.VerticalAlignment = -4108
$ws.Rows.Item(1).RowHeight = 14.4

$ws.Columns.Item(2).ColumnWidth = 10.2222222222222
$ws.Columns.Item(3).ColumnWidth = 9.44444444444444
$ws.Columns.Item(4).ColumnWidth = 12.8888888888889
$ws.Columns.Item(6).ColumnWidth = 11.2222222222222

$ws.Range("E11").Select()
